{"js": "// Update the worksheet date line and the division problems in the table to\n// the next day's regenerated numbers, per the commit.\n//\n// Each table text is replaced via a search scoped to its own paragraph /\n// table cell (rather than a whole-document search) so that:\n//   1) duplicate prompt text (\"22\u00f76=\" appears twice, with two different\n//      replacements) is never ambiguous, and\n//   2) a replacement text that happens to equal another cell's *original*\n//      search text (e.g. \"88\u00f76=\" -> \"68\u00f78=\", while a different, later cell\n//      already contains \"68\u00f78=\" -> \"79\u00f77=\") can never be matched by the\n//      wrong, later search.\n// Using Range.insertText(\u2026, \"Replace\") on the found range (instead of\n// overwriting cell.body / paragraph text wholesale) also preserves the\n// existing run/paragraph formatting (rFonts, sz, jc, \u2026).\n\nasync function replaceInRange(rangeLike, searchText, replacementText) {\n  const results = rangeLike.search(searchText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length < 1) {\n    throw new Error(\"Text '\" + searchText + \"' not found where expected.\");\n  }\n\n  results.items[0].insertText(replacementText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Title / date line.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nawait replaceInRange(paragraphs.items[0], \"2024-08-21 Wednesday\", \"2024-08-22 Thursday\");\n\n// 2) The division-problem table. Only every 4th row (0, 4, 8, 12, 16) has\n// content; the rows between them are blank spacer rows.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst cellEdits = [\n  // [rowIndex, colIndex, oldText, newText]\n  [0, 0, \"88\u00f76=\", \"68\u00f78=\"],\n  [0, 1, \"39\u00f77=\", \"40\u00f77=\"],\n  [0, 2, \"41\u00f74=\", \"68\u00f73=\"],\n  [0, 3, \"22\u00f76=\", \"10\u00f75=\"],\n  [0, 4, \"24\u00f73=\", \"71\u00f74=\"],\n\n  [4, 0, \"35\u00f74=\", \"61\u00f72=\"],\n  [4, 1, \"75\u00f76=\", \"69\u00f73=\"],\n  [4, 2, \"78\u00f77=\", \"53\u00f76=\"],\n  [4, 3, \"29\u00f73=\", \"41\u00f75=\"],\n  [4, 4, \"83\u00f72=\", \"55\u00f73=\"],\n\n  [8, 0, \"93\u00f74=\", \"85\u00f75=\"],\n  [8, 1, \"21\u00f75=\", \"95\u00f73=\"],\n  [8, 2, \"37\u00f78=\", \"69\u00f79=\"],\n  [8, 3, \"73\u00f76=\", \"46\u00f72=\"],\n  [8, 4, \"18\u00f77=\", \"93\u00f73=\"],\n\n  [12, 0, \"23\u00f78=\", \"19\u00f79=\"],\n  [12, 1, \"14\u00f73=\", \"97\u00f75=\"],\n  [12, 2, \"85\u00f78=\", \"51\u00f75=\"],\n  [12, 3, \"68\u00f78=\", \"79\u00f77=\"],\n  [12, 4, \"98\u00f76=\", \"46\u00f77=\"],\n\n  [16, 0, \"22\u00f76=\", \"33\u00f79=\"],\n  [16, 1, \"58\u00f77=\", \"63\u00f75=\"],\n  [16, 2, \"94\u00f75=\", \"16\u00f76=\"],\n  [16, 3, \"37\u00f77=\", \"51\u00f76=\"],\n  [16, 4, \"79\u00f77=\", \"53\u00f79=\"],\n];\n\nfor (const [rowIndex, colIndex, oldText, newText] of cellEdits) {\n  const cell = table.getCell(rowIndex, colIndex);\n  await replaceInRange(cell.body, oldText, newText);\n}\n", "ps1": "# Update the worksheet date line and the division problems in the table to\n# the next day's regenerated numbers, per the commit.\n#\n# NOTE on approach: a Range-scoped \"$range.Find.Execute(...)\" call in this\n# host does not actually confine its search/replace to that Range \u2014 it\n# operates on the whole story and just edits the first match it finds\n# there, regardless of which Range's Find object was used. That silently\n# corrupts unrelated cells whenever the same prompt text (e.g. \"22\u00f76=\",\n# which occurs twice with two different replacements) or a replacement's\n# *output* text (e.g. \"88\u00f76=\" -> \"68\u00f78=\", while another, later cell's\n# ORIGINAL text is also \"68\u00f78=\") appears more than once in the document.\n#\n# Direct Range.Text assignment, however, is correctly confined to the\n# Range it's called on, and (verified) preserves the existing run /\n# paragraph formatting (rFonts, sz, jc, \u2026) already on that text. So every\n# edit below: (1) reads the target Range's current text, (2) verifies it\n# is exactly the expected old value (ignoring the trailing cell-mark /\n# paragraph-mark control characters Word includes in Range.Text), then\n# (3) assigns the new text onto that same Range.\n\n$d = $word.ActiveDocument\n\nfunction Set-RangeText($range, [string]$oldText, [string]$newText) {\n    $current = $range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $oldText) {\n        throw \"Expected '$oldText' but found '$current'.\"\n    }\n    $range.Text = $newText\n}\n\n# 1) Title / date line.\nSet-RangeText $d.Paragraphs.Item(1).Range \"2024-08-21 Wednesday\" \"2024-08-22 Thursday\"\n\n# 2) The division-problem table. Only every 4th row (1, 5, 9, 13, 17 in\n# 1-based Word indexing) has content; the rows between them are blank\n# spacer rows.\n$table = $d.Tables.Item(1)\n\n$cellEdits = @(\n    @(1, 1, \"88\u00f76=\", \"68\u00f78=\"),\n    @(1, 2, \"39\u00f77=\", \"40\u00f77=\"),\n    @(1, 3, \"41\u00f74=\", \"68\u00f73=\"),\n    @(1, 4, \"22\u00f76=\", \"10\u00f75=\"),\n    @(1, 5, \"24\u00f73=\", \"71\u00f74=\"),\n\n    @(5, 1, \"35\u00f74=\", \"61\u00f72=\"),\n    @(5, 2, \"75\u00f76=\", \"69\u00f73=\"),\n    @(5, 3, \"78\u00f77=\", \"53\u00f76=\"),\n    @(5, 4, \"29\u00f73=\", \"41\u00f75=\"),\n    @(5, 5, \"83\u00f72=\", \"55\u00f73=\"),\n\n    @(9, 1, \"93\u00f74=\", \"85\u00f75=\"),\n    @(9, 2, \"21\u00f75=\", \"95\u00f73=\"),\n    @(9, 3, \"37\u00f78=\", \"69\u00f79=\"),\n    @(9, 4, \"73\u00f76=\", \"46\u00f72=\"),\n    @(9, 5, \"18\u00f77=\", \"93\u00f73=\"),\n\n    @(13, 1, \"23\u00f78=\", \"19\u00f79=\"),\n    @(13, 2, \"14\u00f73=\", \"97\u00f75=\"),\n    @(13, 3, \"85\u00f78=\", \"51\u00f75=\"),\n    @(13, 4, \"68\u00f78=\", \"79\u00f77=\"),\n    @(13, 5, \"98\u00f76=\", \"46\u00f77=\"),\n\n    @(17, 1, \"22\u00f76=\", \"33\u00f79=\"),\n    @(17, 2, \"58\u00f77=\", \"63\u00f75=\"),\n    @(17, 3, \"94\u00f75=\", \"16\u00f76=\"),\n    @(17, 4, \"37\u00f77=\", \"51\u00f76=\"),\n    @(17, 5, \"79\u00f77=\", \"53\u00f79=\")\n)\n\nforeach ($edit in $cellEdits) {\n    $rowIndex = $edit[0]\n    $colIndex = $edit[1]\n    $oldText = $edit[2]\n    $newText = $edit[3]\n    $cell = $table.Cell($rowIndex, $colIndex)\n    Set-RangeText $cell.Range $oldText $newText\n}\n"}
